# The "Ölfilter" / "Louis" row (row 7) was removed from the "Tabelle1" table
# on the worksheet. All rows below it shift up by one; the Excel Table
# (ListObject), its AutoFilter/sortState ranges and the "Preis" totals-row
# formula all auto-adjust to the new, smaller range as a consequence.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the sheet row (it belongs to the table, so the table/autofilter/
# totals formula/sortState all shrink from A1:E9 to A1:E8 automatically).
$ws.Rows.Item(7).Delete()

# Column D ("Preis") no longer carries the extra numeric-format style that
# used to be applied to it (cellXfs shrinks back down since that xf becomes
# unused) - reset it back to the plain default style without touching any
# other rows/columns.
$ws.Range("D1:D8").Style = "Standard"
$ws.Columns.Item(4).ClearFormats()
$ws.Range("D9:D11").Clear()

# Leave the selection where the workbook was last saved.
$ws.Range("B13").Select() | Out-Null
